$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 1114. This shifts the existing rows 1114-1162
# down to 1117-1165 (Excel carries every cell value/format along), and the
# former last group (old rows 1160-1162) lands at the new end (1163-1165)
# automatically since nothing is deleted - the sheet simply grows by 3 rows.
$ws.Range("A1114:A1116").EntireRow.Insert()

# Populate the 3 newly-inserted blank rows with the new week's data
# (fecha 2023-05-29 = serial 45075), keeping the same constant columns
# (A,B,C,E-K,Q,R,T) as the rest of this market/product block.

# Row 1114 - Especial
$ws.Range("A1114").Value2 = 8
$ws.Range("B1114").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1114").Value = "Coquimbo"
$ws.Range("D1114").Value2 = 45075
$ws.Range("E1114").Value2 = 4
$ws.Range("F1114").Value = "Fruta"
$ws.Range("G1114").Value2 = 100101
$ws.Range("H1114").Value = "Berries"
$ws.Range("I1114").Value2 = 100112025
$ws.Range("J1114").Value = "Frutilla"
$ws.Range("K1114").Value = "Sin especificar"
$ws.Range("L1114").Value = "Especial"
$ws.Range("M1114").Value2 = 400
$ws.Range("N1114").Value2 = 20000
$ws.Range("O1114").Value2 = 21000
$ws.Range("P1114").Value2 = 20500
$ws.Range("Q1114").Value = "`$/bandeja 7 kilos"
$ws.Range("R1114").Value = "Provincia de Melipilla"
$ws.Range("S1114").Value2 = 2929
$ws.Range("T1114").Value2 = 7

# Row 1115 - Primera
$ws.Range("A1115").Value2 = 8
$ws.Range("B1115").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1115").Value = "Coquimbo"
$ws.Range("D1115").Value2 = 45075
$ws.Range("E1115").Value2 = 4
$ws.Range("F1115").Value = "Fruta"
$ws.Range("G1115").Value2 = 100101
$ws.Range("H1115").Value = "Berries"
$ws.Range("I1115").Value2 = 100112025
$ws.Range("J1115").Value = "Frutilla"
$ws.Range("K1115").Value = "Sin especificar"
$ws.Range("L1115").Value = "Primera"
$ws.Range("M1115").Value2 = 500
$ws.Range("N1115").Value2 = 17000
$ws.Range("O1115").Value2 = 18000
$ws.Range("P1115").Value2 = 17500
$ws.Range("Q1115").Value = "`$/bandeja 7 kilos"
$ws.Range("R1115").Value = "Provincia de Melipilla"
$ws.Range("S1115").Value2 = 2500
$ws.Range("T1115").Value2 = 7

# Row 1116 - Segunda
$ws.Range("A1116").Value2 = 8
$ws.Range("B1116").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1116").Value = "Coquimbo"
$ws.Range("D1116").Value2 = 45075
$ws.Range("E1116").Value2 = 4
$ws.Range("F1116").Value = "Fruta"
$ws.Range("G1116").Value2 = 100101
$ws.Range("H1116").Value = "Berries"
$ws.Range("I1116").Value2 = 100112025
$ws.Range("J1116").Value = "Frutilla"
$ws.Range("K1116").Value = "Sin especificar"
$ws.Range("L1116").Value = "Segunda"
$ws.Range("M1116").Value2 = 360
$ws.Range("N1116").Value2 = 13000
$ws.Range("O1116").Value2 = 14000
$ws.Range("P1116").Value2 = 13500
$ws.Range("Q1116").Value = "`$/bandeja 7 kilos"
$ws.Range("R1116").Value = "Provincia de Melipilla"
$ws.Range("S1116").Value2 = 1929
$ws.Range("T1116").Value2 = 7

Write-Output "Done"
